$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1) - insert "kurtosis" after "skewness", remove "spread"
$headers = @(
    "log_ret",
    "volatility",
    "skewness",
    "kurtosis",
    "amihud",
    "maxmin_ratio",
    "btc_volume",
    "time",
    "delta_5",
    "vol_pre",
    "open_interest",
    "slope",
    "volume",
    "contract_is_call",
    "inter_call_money",
    "inter_put_money",
    "inter_call_skewness"
)
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Update correlation matrix values (rows 2-18, cols A-Q)
$data = @(
    @(1, -0.02596062146849555, -0.06792257891305753, 0.127374286705381, 0.05321439153365063, -0.08318205529627655, 0.0006222378705329488, 0.06599200862823393, 0.103534524432653, -0.03230258034257761, 0.0004237881750733119, 0.01709731236874505, 0.02940419783047241, 0.06818293347164253, 0.06529255357780048, -0.05751957549143164, -0.08590055643659278),
    @(-0.02596062146849555, 1, 0.5542736591930294, -0.6722939266465148, 0.03386044397493845, -0.2953946030403048, 0.05499020499409717, 0.02790638201429849, -0.091105391136891, 0.1302176689297282, 0.1707802744286208, -0.1532407111786153, 0.1458546538442759, -0.03311648640182245, -0.08003899793114443, 0.00807365817573868, 0.4504485234429239),
    @(-0.06792257891305753, 0.5542736591930294, 1, -0.1437616467305882, -0.3062468081764679, -0.04166653077789922, 0.05489822547336181, 0.009047622120908805, 0.0001839619085308236, -0.02727926313903536, 0.05611398453317493, -0.07708101245672777, 0.08132041989108033, 0.03899178051332702, 0.03341198970895573, -0.05726915727237168, 0.776460853462924),
    @(0.127374286705381, -0.6722939266465148, -0.1437616467305882, 1, 0.1856823079631591, 0.4981285502204547, 0.2419691959994472, -0.01840164647348858, 0.07047187697880532, 0.01543339531675838, -0.1684992734965898, 0.1451000279410709, -0.1477892538188408, 0.005541542908669071, 0.07435966035994092, 0.01523555679275454, -0.1714106362870355),
    @(0.05321439153365063, 0.03386044397493845, -0.3062468081764679, 0.1856823079631591, 1, 0.6361530941453635, 0.7273524752812367, 0.192482637529423, -0.1243879243073078, 0.534804614702407, -0.1620391039947083, 0.09664407140331076, -0.1757230374760637, -0.1967552006936234, -0.1563417380003906, 0.2087422891114668, -0.2416223378585217),
    @(-0.08318205529627655, -0.2953946030403048, -0.04166653077789922, 0.4981285502204547, 0.6361530941453635, 1, 0.6960005294985676, 0.1513439541284659, -0.06565636359159345, 0.3349955595730471, -0.2235866524158559, 0.09111724659046457, -0.2040893177355899, -0.1275519774768124, -0.06967366847444846, 0.1333274923959373, -0.0245743112787816),
    @(0.0006222378705329488, 0.05499020499409717, 0.05489822547336181, 0.2419691959994472, 0.7273524752812367, 0.6960005294985676, 1, 0.1661275482465582, -0.07107754698627147, 0.4718195929981553, -0.1485959568250405, 0.1427363978668229, -0.1453291580117486, -0.1492258848970509, -0.08916160327741866, 0.1664096994858938, 0.02492257974378143),
    @(0.06599200862823393, 0.02790638201429849, 0.009047622120908805, -0.01840164647348858, 0.192482637529423, 0.1513439541284659, 0.1661275482465582, 1, 0.06114989471099779, 0.3317632867883254, -0.08169209694182135, -0.07902730685905968, -0.09386926640037982, -0.03838801205346403, -0.006457698279505384, 0.02248182617776148, 0.005003817325760232),
    @(0.103534524432653, -0.091105391136891, 0.0001839619085308236, 0.07047187697880532, -0.1243879243073078, -0.06565636359159345, -0.07107754698627147, 0.06114989471099779, 1, -0.147619821283394, 0.1253831701299011, 0.05149792729412004, 0.1548473616074248, 0.9412385020491509, 0.9435525692865563, -0.8884624562133482, -0.1196685263880993),
    @(-0.03230258034257761, 0.1302176689297282, -0.02727926313903536, 0.01543339531675838, 0.534804614702407, 0.3349955595730471, 0.4718195929981553, 0.3317632867883254, -0.147619821283394, 1, -0.07572380681123662, -0.0865813881341323, -0.1183019606698989, -0.2105239057076246, -0.1534732785628016, 0.1962347078690383, -0.05377120925591405),
    @(0.0004237881750733119, 0.1707802744286208, 0.05611398453317493, -0.1684992734965898, -0.1620391039947083, -0.2235866524158559, -0.1485959568250405, -0.08169209694182135, 0.1253831701299011, -0.07572380681123662, 1, -0.1217153524375737, 0.3418468343867793, 0.1895203005331416, 0.1391149641212611, -0.1914223894951326, 0.01963526155135424),
    @(0.01709731236874505, -0.1532407111786153, -0.07708101245672777, 0.1451000279410709, 0.09664407140331076, 0.09111724659046457, 0.1427363978668229, -0.07902730685905968, 0.05149792729412004, -0.0865813881341323, -0.1217153524375737, 1, -0.07560558277854092, -0.1901587693476353, -0.1194581433103747, 0.3231529086341685, 0.02214667692113197),
    @(0.02940419783047241, 0.1458546538442759, 0.08132041989108033, -0.1477892538188408, -0.1757230374760637, -0.2040893177355899, -0.1453291580117486, -0.09386926640037982, 0.1548473616074248, -0.1183019606698989, 0.3418468343867793, -0.07560558277854092, 1, 0.199060975190523, 0.1608775445991989, -0.1925875752328105, 0.06964651573374943),
    @(0.06818293347164253, -0.03311648640182245, 0.03899178051332702, 0.005541542908669071, -0.1967552006936234, -0.1275519774768124, -0.1492258848970509, -0.03838801205346403, 0.9412385020491509, -0.2105239057076246, 0.1895203005331416, -0.1901587693476353, 0.199060975190523, 1, 0.9488590992266038, -0.9828670152480208, -0.1058421727562293),
    @(0.06529255357780048, -0.08003899793114443, 0.03341198970895573, 0.07435966035994092, -0.1563417380003906, -0.06967366847444846, -0.08916160327741866, -0.006457698279505384, 0.9435525692865563, -0.1534732785628016, 0.1391149641212611, -0.1194581433103747, 0.1608775445991989, 0.9488590992266038, 1, -0.9326023107477749, -0.1049714584497847),
    @(-0.05751957549143164, 0.00807365817573868, -0.05726915727237168, 0.01523555679275454, 0.2087422891114668, 0.1333274923959373, 0.1664096994858938, 0.02248182617776148, -0.8884624562133482, 0.1962347078690383, -0.1914223894951326, 0.3231529086341685, -0.1925875752328105, -0.9828670152480208, -0.9326023107477749, 1, 0.1040287804242795),
    @(-0.08590055643659278, 0.4504485234429239, 0.776460853462924, -0.1714106362870355, -0.2416223378585217, -0.0245743112787816, 0.02492257974378143, 0.005003817325760232, -0.1196685263880993, -0.05377120925591405, 0.01963526155135424, 0.02214667692113197, 0.06964651573374943, -0.1058421727562293, -0.1049714584497847, 0.1040287804242795, 1)
)
for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 1).Value = $row[$c]
    }
}